$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 103
$ws.Range("H103").Value = 5464045.5
$ws.Range("I103").Value = 15025226
$ws.Range("J103").Value = 514
$ws.Range("K103").Value = 45075678
$ws.Range("L103").Value = 1542
$ws.Range("M103").Value = -45075092
$ws.Range("N103").Value = -2714

# Row 107
$ws.Range("H107").Value = 8924.857
$ws.Range("I107").Value = 12209.9
$ws.Range("J107").Value = 712.25
$ws.Range("K107").Value = 12209.9
$ws.Range("L107").Value = 712.25
$ws.Range("M107").Value = -10289.9
$ws.Range("N107").Value = -4552.25

# Row 137
$ws.Range("H137").Value = 7693211
$ws.Range("I137").Value = 855.9231
$ws.Range("J137").Value = 15385566
$ws.Range("K137").Value = 2567.7693
$ws.Range("L137").Value = 46156698
$ws.Range("M137").Value = -17.76929999999993
$ws.Range("N137").Value = -46161798

# Row 138
$ws.Range("H138").Value = 1867.9056
$ws.Range("I138").Value = 1099.9756
$ws.Range("J138").Value = 4491.6665
$ws.Range("K138").Value = 3299.9268
$ws.Range("L138").Value = 13474.9995
$ws.Range("M138").Value = 1840.0732
$ws.Range("N138").Value = -23754.9995

# Row 141
$ws.Range("H141").Value = 1127.5491
$ws.Range("I141").Value = 828.7778
$ws.Range("J141").Value = 3368.3333
$ws.Range("K141").Value = 2486.3334
$ws.Range("L141").Value = 10104.9999
$ws.Range("M141").Value = 2693.6666
$ws.Range("N141").Value = -20464.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1235.1428
$ws.Range("I2").Value = 744.6
$ws.Range("J2").Value = 1507.6666
$ws.Range("K2").Value = 744.6
$ws.Range("L2").Value = 1507.6666
$ws.Range("M2").Value = -631.6
$ws.Range("N2").Value = -1733.6666

# Row 32
$ws.Range("H32").Value = 9512.652
$ws.Range("I32").Value = 9314.161
$ws.Range("J32").Value = 10367.692
$ws.Range("K32").Value = 9314.161
$ws.Range("L32").Value = 10367.692
$ws.Range("M32").Value = -9027.161
$ws.Range("N32").Value = -10941.692

# Row 61
$ws.Range("H61").Value = 10417989
$ws.Range("I61").Value = 13159155
$ws.Range("J61").Value = 1557.9
$ws.Range("K61").Value = 13159155
$ws.Range("L61").Value = 1557.9
$ws.Range("M61").Value = -13158943
$ws.Range("N61").Value = -1981.9

# Row 110
$ws.Range("H110").Value = 1455.8
$ws.Range("I110").Value = 1207.3572
$ws.Range("J110").Value = 2035.5
$ws.Range("K110").Value = 1207.3572
$ws.Range("L110").Value = 2035.5
$ws.Range("M110").Value = 837.6428000000001
$ws.Range("N110").Value = -6125.5

# Row 116
$ws.Range("H116").Value = 1235.1428
$ws.Range("I116").Value = 744.6
$ws.Range("J116").Value = 1507.6666
$ws.Range("K116").Value = 744.6
$ws.Range("L116").Value = 1507.6666
$ws.Range("M116").Value = 1549.4
$ws.Range("N116").Value = -6095.6666

# Row 122
$ws.Range("H122").Value = 7275.3887
$ws.Range("I122").Value = 7684.8125
$ws.Range("K122").Value = 23054.4375
$ws.Range("M122").Value = -20604.4375

# Row 132
$ws.Range("H132").Value = 5954415.5
$ws.Range("I132").Value = 8335033
$ws.Range("J132").Value = 2872.5
$ws.Range("K132").Value = 25005099
$ws.Range("L132").Value = 8617.5
$ws.Range("M132").Value = -25002569
$ws.Range("N132").Value = -13677.5

# Row 136
$ws.Range("H136").Value = 10417989
$ws.Range("I136").Value = 13159155
$ws.Range("J136").Value = 1557.9
$ws.Range("K136").Value = 39477465
$ws.Range("L136").Value = 4673.700000000001
$ws.Range("M136").Value = -39474915
$ws.Range("N136").Value = -9773.700000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1235.1428
$ws.Range("I3").Value = 744.6
$ws.Range("J3").Value = 1507.6666
$ws.Range("K3").Value = 744.6
$ws.Range("L3").Value = 1507.6666
$ws.Range("M3").Value = -630.6
$ws.Range("N3").Value = -1735.6666

# Row 107
$ws.Range("H107").Value = 1241.3846
$ws.Range("I107").Value = 1141.7142
$ws.Range("J107").Value = 1660
$ws.Range("K107").Value = 1141.7142
$ws.Range("L107").Value = 1660
$ws.Range("M107").Value = 778.2858000000001
$ws.Range("N107").Value = -5500

# Row 130
$ws.Range("H130").Value = 63900
$ws.Range("J130").Value = 63900
$ws.Range("L130").Value = 63900
$ws.Range("N130").Value = -73940

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 978.1070999999999
$ws.Range("I16").Value = 979.7
$ws.Range("J16").Value = 974.125
$ws.Range("K16").Value = 979.7
$ws.Range("L16").Value = 974.125
$ws.Range("M16").Value = -692.7
$ws.Range("N16").Value = -1548.125

# Row 31
$ws.Range("H31").Value = 6670643.5
$ws.Range("I31").Value = 4856.788
$ws.Range("J31").Value = 19610112
$ws.Range("K31").Value = 4856.788
$ws.Range("L31").Value = 19610112
$ws.Range("M31").Value = -4561.788
$ws.Range("N31").Value = -19610702

# Row 34
$ws.Range("H34").Value = 6670643.5
$ws.Range("I34").Value = 4856.788
$ws.Range("J34").Value = 19610112
$ws.Range("K34").Value = 4856.788
$ws.Range("L34").Value = 19610112
$ws.Range("M34").Value = -4654.788
$ws.Range("N34").Value = -19610516

# Row 41
$ws.Range("H41").Value = 21032.5
$ws.Range("J41").Value = 30065
$ws.Range("L41").Value = 30065
$ws.Range("N41").Value = -30921

# Row 50
$ws.Range("H50").Value = 9492
$ws.Range("J50").Value = 9492
$ws.Range("L50").Value = 9492
$ws.Range("N50").Value = -10742

# Row 51
$ws.Range("H51").Value = 9531.666999999999
$ws.Range("J51").Value = 9531.666999999999
$ws.Range("L51").Value = 9531.666999999999
$ws.Range("N51").Value = -11003.667

# Row 59
$ws.Range("H59").Value = 15995.25
$ws.Range("J59").Value = 15995.25
$ws.Range("L59").Value = 15995.25
$ws.Range("N59").Value = -18285.25

# Row 60
$ws.Range("H60").Value = 9223
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 8964
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 8964
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -9986

# Row 61
$ws.Range("H61").Value = 9531.666999999999
$ws.Range("J61").Value = 9531.666999999999
$ws.Range("L61").Value = 9531.666999999999
$ws.Range("N61").Value = -10227.667

# Row 113
$ws.Range("H113").Value = 978.1070999999999
$ws.Range("I113").Value = 979.7
$ws.Range("J113").Value = 974.125
$ws.Range("K113").Value = 979.7
$ws.Range("L113").Value = 974.125
$ws.Range("M113").Value = 1190.3
$ws.Range("N113").Value = -5314.125

# Row 122
$ws.Range("H122").Value = 1440.2
$ws.Range("I122").Value = 1440.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4320.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -1870.6

# Row 140
$ws.Range("H140").Value = 46451.285
$ws.Range("J140").Value = 46451.285
$ws.Range("L140").Value = 46451.285
$ws.Range("N140").Value = -56811.285

$ws = $wb.Worksheets.Item("CUL")
# Row 114
$ws.Range("H114").Value = 1622.0555
$ws.Range("J114").Value = 3167.875
$ws.Range("L114").Value = 9503.625
$ws.Range("N114").Value = -16011.625

# Row 118
$ws.Range("H118").Value = 1452.909
$ws.Range("I118").Value = 300
$ws.Range("K118").Value = 900
$ws.Range("M118").Value = 343

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 2606.75
$ws.Range("I107").Value = 2606.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2606.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = -686.75

# Row 113
$ws.Range("H113").Value = 85046.914
$ws.Range("I113").Value = 92623.91
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 92623.91
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = -90453.91
$ws.Range("N113").Value = -6040

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1470.7059
$ws.Range("I61").Value = 1382.0769
$ws.Range("J61").Value = 1758.75
$ws.Range("K61").Value = 1382.0769
$ws.Range("L61").Value = 1758.75
$ws.Range("M61").Value = -1180.0769
$ws.Range("N61").Value = -2162.75

# Row 113
$ws.Range("H113").Value = 1470.7059
$ws.Range("I113").Value = 1382.0769
$ws.Range("J113").Value = 1758.75
$ws.Range("K113").Value = 1382.0769
$ws.Range("L113").Value = 1758.75
$ws.Range("M113").Value = 787.9231
$ws.Range("N113").Value = -6098.75

# Row 122
$ws.Range("H122").Value = 3642.4
$ws.Range("I122").Value = 3593.138
$ws.Range("J122").Value = 3772.2727
$ws.Range("K122").Value = 10779.414
$ws.Range("L122").Value = 11316.8181
$ws.Range("M122").Value = -8329.414000000001
$ws.Range("N122").Value = -16216.8181

# Row 134
$ws.Range("H134").Value = 60513.46
$ws.Range("J134").Value = 60513.46
$ws.Range("L134").Value = 60513.46
$ws.Range("N134").Value = -70653.45999999999
